# Lta-Tnfrsf1a.xlsx -- update scripts wuth new tpm
# Recomputed NATMI TPM-based ligand/receptor edge stats for existing
# sending cluster 'ECs' (rows 2-5), and append rows 6-9 for the newly
# profiled sending cluster 'Resolving-Mac' against the same four
# target clusters (ECs, FAPs, MuSCs, Resolving-Mac).
$rowsData = @(
    @("ECs", "Lta", "Tnfrsf1a", "ECs", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 47.32925566666668, 141.987767, 0.3408416299313156, 0.3408416299313156, 23.63403735775412, 212.706336219787, 0.3266919520370433, 0.3266919520370433),
    @("ECs", "Lta", "Tnfrsf1a", "FAPs", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 43.717953, 131.153859, 0.3148348341399153, 0.3148348341399154, 21.830720129711, 196.476481167399, 0.3017648007232984, 0.3017648007232985),
    @("ECs", "Lta", "Tnfrsf1a", "MuSCs", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 21.09134933333333, 63.274048, 0.1518893501062827, 0.1518893501062827, 10.53204262454756, 94.788383620928, 0.1455838252207007, 0.1455838252207007),
    @("ECs", "Lta", "Tnfrsf1a", "Resolving-Mac", 2, 0.6666666666666666, 0.4993536666666667, 1.498061, 0.9584860631692095, 0.9584860631692095, 3, 1, 26.72140366666666, 80.164211, 0.1924341858224863, 0.1924341858224864, 13.34343089943011, 120.090878094871, 0.184445485188167, 0.1844454851881671),
    @("Resolving-Mac", "Lta", "Tnfrsf1a", "ECs", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 47.32925566666668, 141.987767, 0.3408416299313156, 0.3408416299313156, 1.023637141558667, 9.212734274028001, 0.01414967789427233, 0.01414967789427234),
    @("Resolving-Mac", "Lta", "Tnfrsf1a", "FAPs", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 43.717953, 131.153859, 0.3148348341399153, 0.3148348341399154, 0.945531887484, 8.509786987356, 0.01307003341661687, 0.01307003341661688),
    @("Resolving-Mac", "Lta", "Tnfrsf1a", "MuSCs", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 21.09134933333333, 63.274048, 0.1518893501062827, 0.1518893501062827, 0.4561637033813333, 4.105473330432, 0.006305524885582055, 0.006305524885582058),
    @("Resolving-Mac", "Lta", "Tnfrsf1a", "Resolving-Mac", 1, 0.3333333333333333, 0.021628, 0.064884, 0.04151393683079058, 0.04151393683079059, 3, 1, 26.72140366666666, 80.164211, 0.1924341858224863, 0.1924341858224864, 0.5779305185026665, 5.201374666523999, 0.007988700634319314, 0.007988700634319317)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $rowValues = $rowsData[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
